# Update the dSF column (column F) with repulled/recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -2
    11 = 2
    16 = -1
    22 = -1
    26 = 0
    33 = 2
    34 = 0
    36 = 0
    39 = 2
    43 = -1
    47 = 4
    48 = -2
    52 = 1
    55 = 3
    60 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
